$wb = $excel.ActiveWorkbook

# The workbook originally has a single sheet "Sheet1" holding one metadata
# record (job name, exp folder/replicates, control folder/replicates,
# ionization). We rename it to "All" (keeping its data), relabel the
# "Control ..." headers to "CTRL ...", and add a brand-new first sheet
# "Job to Run" that holds the same record reshaped into the layout used to
# drive the per-run script (job name first, then exp/ctrl folder+replicate
# columns, then ionization).

$allSheet = $wb.ActiveSheet
$allSheet.Name = "All"

# Relabel "Control ..." -> "CTRL ..." on the "All" sheet.
$allSheet.Range("D1").Value = "CTRL Folder Name"
$allSheet.Range("E1").Value = "CTRL num replicates"

# Selection on "All" becomes the data block, and it is no longer the
# visible/active tab.
$allSheet.Range("A1:F2").Select()

# Insert the new "Job to Run" sheet in front of "All" (so it becomes the
# first, left-most tab).
$jobSheet = $wb.Worksheets.Add($allSheet)
$jobSheet.Name = "Job to Run"

# Headers
$jobSheet.Range("A1").Value = "Job Name"
$jobSheet.Range("B1").Value = "Exp Folder Name"
$jobSheet.Range("C1").Value = "EXP num replicates"
$jobSheet.Range("D1").Value = "CTRL Folder Name"
$jobSheet.Range("E1").Value = "CTRL num replicates"
$jobSheet.Range("F1").Value = "Ionization"

# Data row - job to run for TJGIp11 against the TJGI2pt1_EV_gpdA control
$jobSheet.Range("A2").Value = "TJGIp11"
$jobSheet.Range("B2").Value = "Anid_HE_TJGIp11_pos_2018"
$jobSheet.Range("C2").Value = 3
$jobSheet.Range("D2").Value = "Anid_HE_CTRL_TJGI2pt1_EV_gpdA_pos_2018"
$jobSheet.Range("E2").Value = 3
$jobSheet.Range("F2").Value = "POS"

# Match formatting of the data row (vertically centered, as on "All")
$jobSheet.Range("B2:F2").VerticalAlignment = -4108

# Column widths (best-fit sizing as authored; inputs pre-compensated for the
# COM layer's pixel-quantized ColumnWidth rounding so the stored OOXML
# <col> widths land as close as possible to the authored values)
$jobSheet.Columns.Item(1).ColumnWidth = 8.0833333334
$jobSheet.Columns.Item(2).ColumnWidth = 23.5833333334
$jobSheet.Columns.Item(3).ColumnWidth = 15.7500000001
$jobSheet.Columns.Item(4).ColumnWidth = 37.9166666668

$jobSheet.PageSetup.Orientation = 1

$jobSheet.Range("E1").Select()
$jobSheet.Activate()
